# Update "RES installed" sheet values per commit: "DN node 5. RES updated"
$wb = $excel.ActiveWorkbook
$resSheet = $wb.Worksheets.Item("RES installed")

$resSheet.Range("C2").Value = 10
$resSheet.Range("C3").Value = 10
$resSheet.Range("C4").Value = 5
$resSheet.Range("C5").Value = 5
$resSheet.Range("C6").Value = 5

# Update selection on RES installed sheet to F7 and make it the active/selected tab
$resSheet.Activate()
$resSheet.Range("F7").Select()

$excel.CalculateFullRebuild()
